# Natmi following Dr Hou advice
# Update ligand/receptor expressing-cell counts (E, K: 1 -> 3) and the
# resulting average/total expression and specificity values that were
# recomputed downstream (columns G, H, M, N, O, P, Q, R, S, T) for rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, keyed by column letter.
$updates = @{
    2 = @{
        E = 3; G = 4.674406333333334; H = 14.023219; K = 3
        M = 0.502378; N = 1.507134
        O = 0.0189742916423209; P = 0.0189742916423209
        Q = 2.348318904927333; R = 21.134870144346
        S = 0.0189742916423209; T = 0.0189742916423209
    }
    3 = @{
        E = 3; G = 4.674406333333334; H = 14.023219; K = 3
        M = 0.6023626666666666; N = 1.807088
        O = 0.02275060793223323; P = 0.02275060793223323
        Q = 2.815687864030222; R = 25.341190776272
        S = 0.02275060793223323; T = 0.02275060793223323
    }
    4 = @{
        E = 3; G = 4.674406333333334; H = 14.023219; K = 3
        M = 5.474800666666667; N = 16.424402
        O = 0.2067774952981744; P = 0.2067774952981744
        Q = 25.59144291000423; R = 230.322986190038
        S = 0.2067774952981744; T = 0.2067774952981744
    }
    5 = @{
        E = 3; G = 4.674406333333334; H = 14.023219; K = 3
        M = 5.012936333333333; N = 15.038809
        O = 0.1893333624741797; P = 0.1893333624741797
        Q = 23.43250134513011; R = 210.892512106171
        S = 0.1893333624741797; T = 0.1893333624741797
    }
    6 = @{
        E = 3; G = 4.674406333333334; H = 14.023219; K = 3
        M = 7.423250333333333; N = 22.269751
        O = 0.2803684014001858; P = 0.2803684014001858
        Q = 34.69928837205211; R = 312.293595348469
        S = 0.2803684014001858; T = 0.2803684014001858
    }
    7 = @{
        E = 3; G = 4.674406333333334; H = 14.023219; K = 3
        M = 7.461044333333334; N = 22.383133
        O = 0.281795841252906; P = 0.281795841252906
        Q = 34.87595288501412; R = 313.883575965127
        S = 0.281795841252906; T = 0.281795841252906
    }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
